$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from H1 (bold, bordered, centered) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 5
